$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.874.18"
$ws.Range("E2").Value = "  -0.21%  "
$ws.Range("D3").Value = "1.628.07"
$ws.Range("E3").Value = "  -0.81%  "
$ws.Range("D4").Value = "0.997"
$ws.Range("E4").Value = "  -0.28%  "
$ws.Range("D5").Value = "211.55"
$ws.Range("E5").Value = "  -0.67%  "
$ws.Range("E6").Value = "  -0.32%  "
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("D8").Value = "23.36"
$ws.Range("E8").Value = "  -0.49%  "
$ws.Range("E9").Value = "  -0.52%  "
$ws.Range("E10").Value = "  -0.31%  "
$ws.Range("D11").Value = "0.0879"
$ws.Range("E11").Value = "  -0.55%  "
$ws.Range("D12").Value = "1.858.48"
$ws.Range("E12").Value = "  -0.70%  "
$ws.Range("D13").Value = "1.624.51"
$ws.Range("E13").Value = "  -0.94%  "
$ws.Range("E14").Value = "  -1.60%  "
$ws.Range("E15").Value = "  -1.92%  "
$ws.Range("D16").Value = "65.31"
$ws.Range("E16").Value = "  -0.42%  "
$ws.Range("D17").Value = "27.866.17"
$ws.Range("E17").Value = "  -0.21%  "
$ws.Range("D18").Value = "229.95"
$ws.Range("E18").Value = "  -1.31%  "
$ws.Range("D19").Value = "7.64"
$ws.Range("E19").Value = "  +0.55%  "
$ws.Range("D20").Value = "0.0₃0721"
$ws.Range("E20").Value = "  -0.14%  "
$ws.Range("D21").Value = "0.997"
$ws.Range("E21").Value = "  -0.34%  "
$ws.Range("E22").Value = "  -1.26%  "
$ws.Range("D23").Value = "10.09"
$ws.Range("E23").Value = "  -3.80%  "
$ws.Range("D24").Value = "2.04"
$ws.Range("E24").Value = "  -2.86%  "
$ws.Range("D25").Value = "154.32"
$ws.Range("E25").Value = "  +0.96%  "
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("E27").Value = "  -0.35%  "
$ws.Range("D28").Value = "15.50"
$ws.Range("E28").Value = "  -1.40%  "
$ws.Range("D29").Value = "0.997"
$ws.Range("E29").Value = "  -0.32%  "
$ws.Range("D30").Value = "1.18"
$ws.Range("E30").Value = "  -1.13%  "
$ws.Range("D31").Value = "0.0481"
$ws.Range("E31").Value = "  -0.73%  "
$ws.Range("E32").Value = "  +1.68%  "
$ws.Range("E33").Value = "  -0.49%  "
$ws.Range("D34").Value = "1.397.30"
$ws.Range("E34").Value = "  -0.80%  "
$ws.Range("E35").Value = "  +0.13%  "
$ws.Range("D36").Value = "1.02"
$ws.Range("E36").Value = "  +10.05%  "
$ws.Range("E37").Value = "  -1.09%  "
$ws.Range("E38").Value = "  +0.28%  "
$ws.Range("D39").Value = "0.555"
$ws.Range("E39").Value = "  -0.99%  "
$ws.Range("D40").Value = "0.852"
$ws.Range("E40").Value = "  -3.41%  "
$ws.Range("D41").Value = "0.997"
$ws.Range("E41").Value = "  -0.31%  "
$ws.Range("E42").Value = "  -2.19%  "
$ws.Range("D43").Value = "1.85"
$ws.Range("E43").Value = "  -0.10%  "
$ws.Range("D44").Value = "65.67"
$ws.Range("E44").Value = "  -2.75%  "
$ws.Range("E45").Value = "  -1.46%  "
$ws.Range("D46").Value = "1.768.39"
$ws.Range("E46").Value = "  -0.59%  "
$ws.Range("E47").Value = "  -2.87%  "
$ws.Range("D48").Value = "87.99"
$ws.Range("E48").Value = "  +0.18%  "
$ws.Range("E49").Value = "  +1.35%  "
$ws.Range("E50").Value = "  +8.66%  "
$ws.Range("E51").Value = "  -0.65%  "
